$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.853.52'
$ws.Range('E2').Value = '  +1.10%  '
$ws.Range('D3').Value = '1.755.77'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '327.36'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4591'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.55%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3495'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '41.89'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07350'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.082'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.000'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.53'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.972'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.165'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.15%  '
$ws.Range('D16').Value = '1.757.80'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.62'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.12%  '
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06413'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.81'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.747'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('D23').Value = '27.892.60'
$ws.Range('E23').Value = '  +1.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.15'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.157'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +4.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '161.40'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.01'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.75%  '
$ws.Range('D28').Value = '1.960.09'
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.147'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.25'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.066'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.62%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09257'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.665'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.533'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '11.69'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02262'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06085'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2057'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.68%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.891'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6171'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.80%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.178'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.44%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.363'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.755'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.06'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.29%  '
$ws.Range('E45').Value = '  +0.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5780'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '123.67'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.925'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.73%  '
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.121'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.09'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.27%  '
